# #5: property aircraft done
# The "建物" (building) sheet's property_category column (I) was
# mistakenly populated with "land" for every data row; fix it to "building".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

$ws.Range("I2:I6").Value = "building"
